# Updated cryptos list on Thu Feb 23 08:52:07 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures, and
# swaps the HuobiToken/BitcoinCash rows (29/30) to reflect their new order.
#
# Each literal is written as  '<apostrophe><text>'  -- the leading
# apostrophe forces Excel to keep the value as text (matching the
# original inlineStr cell type) instead of re-parsing numeric-looking
# strings such as "1.000" or "0.00001337" into plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.436.46'
$ws.Range("E2").Value = '''  +1.68%  '
$ws.Range("D3").Value = '''1.670.15'
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '''  -0.30%  '
$ws.Range("D5").Value = '''312.15'
$ws.Range("E5").Value = '''  +1.84%  '
$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '''  -0.38%  '
$ws.Range("D7").Value = '''0.3970'
$ws.Range("E7").Value = '''  +1.70%  '
$ws.Range("D8").Value = '''0.3928'
$ws.Range("E8").Value = '''  +2.94%  '
$ws.Range("D9").Value = '''52.37'
$ws.Range("E9").Value = '''  +6.42%  '
$ws.Range("E10").Value = '''  +4.89%  '
$ws.Range("D11").Value = '''0.9995'
$ws.Range("E11").Value = '''  -0.39%  '
$ws.Range("D12").Value = '''0.08572'
$ws.Range("E12").Value = '''  +2.61%  '
$ws.Range("D13").Value = '''24.54'
$ws.Range("E13").Value = '''  +4.34%  '
$ws.Range("D14").Value = '''7.308'
$ws.Range("E14").Value = '''  +4.31%  '
$ws.Range("D15").Value = '''7.944'
$ws.Range("E15").Value = '''  +7.51%  '
$ws.Range("D16").Value = '''0.00001337'
$ws.Range("E16").Value = '''  +5.39%  '
$ws.Range("D17").Value = '''1.663.55'
$ws.Range("E17").Value = '''  +0.25%  '
$ws.Range("D18").Value = '''94.99'
$ws.Range("E18").Value = '''  +0.21%  '
$ws.Range("D19").Value = '''0.07032'
$ws.Range("E19").Value = '''  +2.33%  '
$ws.Range("D20").Value = '''20.63'
$ws.Range("E20").Value = '''  -0.04%  '
$ws.Range("E21").Value = '''  +2.16%  '
$ws.Range("D22").Value = '''0.9999'
$ws.Range("E22").Value = '''  -0.42%  '
$ws.Range("D23").Value = '''13.77'
$ws.Range("E23").Value = '''  +2.36%  '
$ws.Range("D24").Value = '''24.446.30'
$ws.Range("E24").Value = '''  +1.70%  '
$ws.Range("D25").Value = '''2.479'
$ws.Range("E25").Value = '''  +6.83%  '
$ws.Range("D26").Value = '''3.069'
$ws.Range("E26").Value = '''  +15.49%  '
$ws.Range("D27").Value = '''22.55'
$ws.Range("E27").Value = '''  +1.59%  '
$ws.Range("D28").Value = '''157.34'
$ws.Range("E28").Value = '''  +0.03%  '
$ws.Range("B29").Value = '''HuobiToken'
$ws.Range("C29").Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D29").Value = '''5.460'
$ws.Range("E29").Value = '''  +2.94%  '
$ws.Range("B30").Value = '''BitcoinCash'
$ws.Range("C30").Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '''142.68'
$ws.Range("E30").Value = '''  +2.47%  '
$ws.Range("D31").Value = '''8.025'
$ws.Range("E31").Value = '''  -6.88%  '
$ws.Range("D32").Value = '''2.540'
$ws.Range("E32").Value = '''  +5.96%  '
$ws.Range("D33").Value = '''1.847.00'
$ws.Range("E33").Value = '''  +0.54%  '
$ws.Range("D34").Value = '''1.066'
$ws.Range("E34").Value = '''  +13.77%  '
$ws.Range("D35").Value = '''0.03110'
$ws.Range("E35").Value = '''  +8.33%  '
$ws.Range("D36").Value = '''0.08287'
$ws.Range("E36").Value = '''  +4.42%  '
$ws.Range("D37").Value = '''6.930'
$ws.Range("E37").Value = '''  +2.48%  '
$ws.Range("D38").Value = '''11.18'
$ws.Range("E38").Value = '''  +13.80%  '
$ws.Range("D39").Value = '''0.2768'
$ws.Range("E39").Value = '''  +4.21%  '
$ws.Range("D40").Value = '''0.09272'
$ws.Range("D41").Value = '''0.7720'
$ws.Range("E41").Value = '''  +3.55%  '
$ws.Range("D42").Value = '''13.77'
$ws.Range("E42").Value = '''  +6.86%  '
$ws.Range("D43").Value = '''1.447'
$ws.Range("E43").Value = '''  +1.23%  '
$ws.Range("D44").Value = '''16.62'
$ws.Range("E44").Value = '''  +5.53%  '
$ws.Range("D45").Value = '''0.7119'
$ws.Range("E45").Value = '''  +4.63%  '
$ws.Range("D46").Value = '''2.551'
$ws.Range("E46").Value = '''  +4.57%  '
$ws.Range("D47").Value = '''4.128'
$ws.Range("E48").Value = '''  -0.34%  '
$ws.Range("D49").Value = '''0.08450'
$ws.Range("E49").Value = '''  +1.74%  '
$ws.Range("D50").Value = '''137.11'
$ws.Range("E50").Value = '''  +4.53%  '
$ws.Range("D51").Value = '''1.271'
$ws.Range("E51").Value = '''  +2.52%  '
